$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.244.02'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.39%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.269.33'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.83%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.35'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.63%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.526'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.14%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.492'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.22'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.96%  '
$ws.Range("E11").Value = '  -1.72%  '
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("E13").Value = '  +1.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.618.07'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.70'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.34%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.266.07'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.44%  '
$ws.Range("E17").Value = '  -1.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.090.62'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.58%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.28'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.40%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0906'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.02'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.83'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.50%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.35'
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = '  -1.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.97'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.11%  '
$ws.Range("E26").Value = '  -0.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.57'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.26%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.86'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.59'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.39%  '
$ws.Range("E30").Value = '  +0.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '162.54'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.26'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.07%  '
$ws.Range("E34").Value = '  +2.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.75'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.97%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0738'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.27%  '
$ws.Range("E37").Value = '  -0.69%  '
$ws.Range("E38").Value = '  -3.58%  '
$ws.Range("E39").Value = '  -0.95%  '
$ws.Range("E40").Value = '  -1.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.10'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.34'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.953.43'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.10'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0281'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.02'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.58%  '
$ws.Range("E47").Value = '  -2.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.06'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.490.51'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.59%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '92.22'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.33%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.66'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.05%  '
